$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed DQ6 seed location: update the Seed row's address from 02101564 to 020BA3AC
$ws.Range("A2").Value = "020BA3AC"

# Move the active selection, matching the saved view state
$ws.Range("A8").Select()
